# Comprehensive_Etsy_Store_Data.xlsx update
# - Insert a new "Total Items Left in Inventory Item N" column right after each
#   "Reviews Item N" column (N = 1..6).
# - Append 21 new "platform-wide" summary columns after the last existing
#   column (Meta Ad 3 Video Transcript).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-HeaderAfter {
    param($afterColLetter, $headerText)

    $afterCell = $ws.Range($afterColLetter + "1")
    $newCell = $afterCell.Offset(0, 1)

    # Push everything from this column onward one column to the right.
    $newCell.EntireColumn.Insert()

    # Populate the freshly-inserted (now blank) header cell.
    $newCell.Value = $headerText

    # Match the header styling (bold, centered, bordered) used by the rest
    # of row 1 by copying it from the adjacent existing header cell.
    $afterCell.Copy()
    $newCell.PasteSpecial(-4122)  # xlPasteFormats
}

# Work right-to-left so each earlier (more-leftward) insertion point keeps
# referring to the same original column letters.
Insert-HeaderAfter "BG" "Total Items Left in Inventory Item 6"
Insert-HeaderAfter "AY" "Total Items Left in Inventory Item 5"
Insert-HeaderAfter "AQ" "Total Items Left in Inventory Item 4"
Insert-HeaderAfter "AI" "Total Items Left in Inventory Item 3"
Insert-HeaderAfter "AA" "Total Items Left in Inventory Item 2"
Insert-HeaderAfter "S"  "Total Items Left in Inventory Item 1"

# After the six insertions above, the last original column (previously DR1,
# "Meta Ad 3 Video Transcript") has shifted six places right to DX1.
$lastCol = "DX"
$lastCell = $ws.Range($lastCol + "1")

$newHeaders = @(
    "Total Number of Items in Carts (Platform-Wide)",
    "Total Number of Shares (Last 3 Months)",
    "Total Number of Shares (Last 2 Months)",
    "Total Number of Shares (Last 1 Month)",
    "Total Number of Shares (Last 2 Weeks)",
    "Total Number of Shares (Last 3 Days)",
    "Total Number of Impressions (Last 3 Months)",
    "Total Number of Impressions (Last 2 Months)",
    "Total Number of Impressions (Last 1 Month)",
    "Total Number of Impressions (Last 2 Weeks)",
    "Total Number of Impressions (Last 3 Days)",
    "Total Number of Comments (Last 3 Months)",
    "Total Number of Comments (Last 2 Months)",
    "Total Number of Comments (Last 1 Month)",
    "Total Number of Comments (Last 2 Weeks)",
    "Total Number of Comments (Last 3 Days)",
    "Total Number of Likes (Last 3 Months)",
    "Total Number of Likes (Last 2 Months)",
    "Total Number of Likes (Last 1 Month)",
    "Total Number of Likes (Last 2 Weeks)",
    "Total Number of Likes (Last 3 Days)"
)

$cursor = $lastCell
foreach ($headerText in $newHeaders) {
    $cell = $cursor.Offset(0, 1)
    $cell.Value = $headerText
    $lastCell.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $cursor = $cell
}
